# Auto-generated edit script: updates leve market-price columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets to reflect refreshed
# market board data (scheduled runner sync).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 145.66667
$ws.Range("J55").Value = 155.625
$ws.Range("L55").Value = 155.625
$ws.Range("N55").Value = -583.625

$ws.Range("H70").Value = 1315.1364
$ws.Range("I70").Value = 944.875
$ws.Range("J70").Value = 1526.7142
$ws.Range("K70").Value = 2834.625
$ws.Range("L70").Value = 4580.142599999999
$ws.Range("M70").Value = -2564.625
$ws.Range("N70").Value = -5120.142599999999

$ws.Range("H73").Value = 1315.1364
$ws.Range("I73").Value = 944.875
$ws.Range("J73").Value = 1526.7142
$ws.Range("K73").Value = 2834.625
$ws.Range("L73").Value = 4580.142599999999
$ws.Range("M73").Value = -1898.625
$ws.Range("N73").Value = -6452.142599999999

$ws.Range("H121").Value = 872.1429000000001
$ws.Range("J121").Value = 1081
$ws.Range("L121").Value = 3243
$ws.Range("N121").Value = -6737

$ws.Range("H125").Value = 16017182
$ws.Range("I125").Value = 332
$ws.Range("J125").Value = 22423922
$ws.Range("K125").Value = 2988
$ws.Range("L125").Value = 201815298
$ws.Range("M125").Value = -528
$ws.Range("N125").Value = -201820218

$ws.Range("H129").Value = 1080
$ws.Range("I129").Value = 268.2
$ws.Range("J129").Value = 1230.3334
$ws.Range("K129").Value = 804.5999999999999
$ws.Range("L129").Value = 3691.0002
$ws.Range("M129").Value = 4195.4
$ws.Range("N129").Value = -13691.0002

$ws.Range("H138").Value = 9982074
$ws.Range("I138").Value = 3099199.2
$ws.Range("J138").Value = 13516524
$ws.Range("K138").Value = 9297597.600000001
$ws.Range("L138").Value = 40549572
$ws.Range("M138").Value = -9292457.600000001
$ws.Range("N138").Value = -40559852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2361.1853
$ws.Range("I61").Value = 1787.9166
$ws.Range("K61").Value = 1787.9166
$ws.Range("M61").Value = -1575.9166

$ws.Range("H122").Value = 3066.3928
$ws.Range("I122").Value = 2342.95
$ws.Range("K122").Value = 7028.849999999999
$ws.Range("M122").Value = -4578.849999999999

$ws.Range("H133").Value = 47333.332
$ws.Range("J133").Value = 47333.332
$ws.Range("L133").Value = 47333.332
$ws.Range("N133").Value = -52393.332

$ws.Range("H136").Value = 2361.1853
$ws.Range("I136").Value = 1787.9166
$ws.Range("K136").Value = 5363.7498
$ws.Range("M136").Value = -2813.7498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 469.25
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 638.5
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 638.5
$ws.Range("M22").Value = -127
$ws.Range("N22").Value = -984.5

$ws.Range("H94").Value = 1276.44
$ws.Range("I94").Value = 1264.1364
$ws.Range("J94").Value = 1366.6666
$ws.Range("K94").Value = 1264.1364
$ws.Range("L94").Value = 1366.6666
$ws.Range("M94").Value = -813.1364000000001
$ws.Range("N94").Value = -2268.6666

$ws.Range("H134").Value = 3390.068
$ws.Range("I134").Value = 1856.8125
$ws.Range("J134").Value = 7478.75
$ws.Range("K134").Value = 5570.4375
$ws.Range("L134").Value = 22436.25
$ws.Range("M134").Value = -3035.4375
$ws.Range("N134").Value = -27506.25

$ws.Range("H140").Value = 64800
$ws.Range("J140").Value = 64800
$ws.Range("L140").Value = 64800
$ws.Range("N140").Value = -75160

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6278.225
$ws.Range("I31").Value = 2050.6155
$ws.Range("K31").Value = 2050.6155
$ws.Range("M31").Value = -1755.6155

$ws.Range("H34").Value = 6278.225
$ws.Range("I34").Value = 2050.6155
$ws.Range("K34").Value = 2050.6155
$ws.Range("M34").Value = -1848.6155

$ws.Range("H58").Value = 2116.5813
$ws.Range("I58").Value = 1126.4333
$ws.Range("J58").Value = 4401.5386
$ws.Range("K58").Value = 1126.4333
$ws.Range("L58").Value = 4401.5386
$ws.Range("M58").Value = -923.4332999999999
$ws.Range("N58").Value = -4807.5386

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H134").Value = 3197.875
$ws.Range("I134").Value = 1545.6
$ws.Range("J134").Value = 5951.6665
$ws.Range("K134").Value = 4636.799999999999
$ws.Range("L134").Value = 17854.9995
$ws.Range("M134").Value = -2101.799999999999
$ws.Range("N134").Value = -22924.9995

$ws.Range("H136").Value = 2116.5813
$ws.Range("I136").Value = 1126.4333
$ws.Range("J136").Value = 4401.5386
$ws.Range("K136").Value = 3379.2999
$ws.Range("L136").Value = 13204.6158
$ws.Range("M136").Value = -829.2999
$ws.Range("N136").Value = -18304.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 155.62962
$ws.Range("I38").Value = 150
$ws.Range("J38").Value = 158
$ws.Range("K38").Value = 450
$ws.Range("L38").Value = 474
$ws.Range("M38").Value = -103
$ws.Range("N38").Value = -1168

$ws.Range("H113").Value = 561.9
$ws.Range("I113").Value = 460
$ws.Range("J113").Value = 573.2222
$ws.Range("K113").Value = 1380
$ws.Range("L113").Value = 1719.6666
$ws.Range("M113").Value = 790
$ws.Range("N113").Value = -6059.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4605.636
$ws.Range("I70").Value = 4688.9443
$ws.Range("K70").Value = 4688.9443
$ws.Range("M70").Value = -4418.9443

$ws.Range("H73").Value = 4605.636
$ws.Range("I73").Value = 4688.9443
$ws.Range("K73").Value = 4688.9443
$ws.Range("M73").Value = -3752.9443

$ws.Range("H126").Value = 2478.1482
$ws.Range("I126").Value = 1978.2222
$ws.Range("J126").Value = 2728.111
$ws.Range("K126").Value = 5934.6666
$ws.Range("L126").Value = 8184.333
$ws.Range("M126").Value = -3464.6666
$ws.Range("N126").Value = -13124.333

$ws.Range("H137").Value = 44560
$ws.Range("J137").Value = 44560
$ws.Range("L137").Value = 44560
$ws.Range("N137").Value = -54760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3064.95
$ws.Range("I7").Value = 2200
$ws.Range("J7").Value = 3217.5881
$ws.Range("K7").Value = 2200
$ws.Range("L7").Value = 3217.5881
$ws.Range("M7").Value = -2088
$ws.Range("N7").Value = -3441.5881

$ws.Range("H40").Value = 2787.48
$ws.Range("I40").Value = 2015.6666
$ws.Range("J40").Value = 3499.923
$ws.Range("K40").Value = 2015.6666
$ws.Range("L40").Value = 3499.923
$ws.Range("M40").Value = -1879.6666
$ws.Range("N40").Value = -3771.923

$ws.Range("H126").Value = 3064.95
$ws.Range("I126").Value = 2200
$ws.Range("J126").Value = 3217.5881
$ws.Range("K126").Value = 6600
$ws.Range("L126").Value = 9652.764299999999
$ws.Range("M126").Value = -4130
$ws.Range("N126").Value = -14592.7643

$ws.Range("H132").Value = 3014.9768
$ws.Range("I132").Value = 1805.7241
$ws.Range("J132").Value = 5519.857
$ws.Range("K132").Value = 5417.1723
$ws.Range("L132").Value = 16559.571
$ws.Range("M132").Value = -2887.1723
$ws.Range("N132").Value = -21619.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19383.389
$ws.Range("J62").Value = 16766.555
$ws.Range("L62").Value = 16766.555
$ws.Range("N62").Value = -18014.555

$ws.Range("H65").Value = 19383.389
$ws.Range("J65").Value = 16766.555
$ws.Range("L65").Value = 83832.77499999999
$ws.Range("N65").Value = -90072.77499999999

$ws.Range("H132").Value = 3291.4878
$ws.Range("I132").Value = 3352.125
$ws.Range("J132").Value = 3075.889
$ws.Range("K132").Value = 10056.375
$ws.Range("L132").Value = 9227.667000000001
$ws.Range("M132").Value = -7526.375
$ws.Range("N132").Value = -14287.667
